# Thicken the outline weight of the logo/icon shapes on the last slide.
#
# Before:
#   Rectangle 75  (no fill)            a:ln w="3175"  (0.25 pt)
#   Rectangle 148 (fill 2C80B8)        a:ln w="3175"  (0.25 pt)
#   Rectangle 148 (fill 7FCEBB)        a:ln w="3175"  (0.25 pt)
#   Rectangle 148 (fill EDF9B1)        a:ln w="3175"  (0.25 pt)
#
# After:
#   Rectangle 75  (no fill)            a:ln w="38100" (3 pt)
#   Rectangle 148 (fill 2C80B8)        a:ln w="25400" (2 pt)
#   Rectangle 148 (fill 7FCEBB)        a:ln w="25400" (2 pt)
#   Rectangle 148 (fill EDF9B1)        a:ln w="25400" (2 pt)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# Shape 1: the plain (no-fill) outline rectangle -> 3 pt line
$s.Shapes.Item(1).Line.Weight = 3

# Shapes 2-4: the three colored logo/icon pieces -> 2 pt line
$s.Shapes.Item(2).Line.Weight = 2
$s.Shapes.Item(3).Line.Weight = 2
$s.Shapes.Item(4).Line.Weight = 2
